$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 28055.51662633527
$ws.Range("B2").Value = 0
$ws.Range("J2").Value = 0

$ws.Range("A3").Value = 26905.70154390091
$ws.Range("G3").Value = 0

$ws.Range("A4").Value = 25723.37458066656

$ws.Range("A5").Value = 24845.72510195928

$ws.Range("A6").Value = 24988.90498752365

$ws.Range("A7").Value = 25835.11014617304

$ws.Range("A8").Value = 26401.55040670632

$ws.Range("A9").Value = 28889.35044967622

$ws.Range("A10").Value = 33516.87245290272

$ws.Range("A11").Value = 29249.825

$ws.Range("A12").Value = 29168.172
$ws.Range("B12").Value = 10000
$ws.Range("J12").Value = 10000

$ws.Range("A13").Value = 30065.795
$ws.Range("B13").Value = 10000
$ws.Range("J13").Value = 10000

$ws.Range("A14").Value = 30878.246

$ws.Range("A15").Value = 32676.0655

$ws.Range("A16").Value = 49763.3545

$ws.Range("A17").Value = 28812.386

$ws.Range("A18").Value = 96.81750000000102

$ws.Range("C23").Value = 24000

$ws.Range("C24").Value = 21000

$ws.Range("C25").Value = 17500
